# "changed how to read datasets"
#
# The header labels are renamed, and the three "Tasso di cambio medio"
# sample values for the EUR/USD and EUR/JPY "BUDGET"/"CONSUNTIVO" rows
# (C3, C6, C7) were re-entered using a comma decimal separator, which
# Excel stored as literal text instead of numbers (C2/C4/C5 keep their
# original numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column headers.
$ws.Range("A1").Value = "Codice valuta"
$ws.Range("C1").Value = "Tasso di cambio medio"

# Re-enter the exchange-rate figures using a comma decimal separator so
# they land in the sheet as text (matching the new values verbatim).
$ws.Range("C3").Value = "1,0541"

$ws.Range("C6").Value = "1,1993"
$ws.Range("C6").NumberFormat = "0"

$ws.Range("C7").Value = "135,01"
$ws.Range("C7").NumberFormat = "0"
